$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = " Abu Dhabi"
$ws.Range("B4").Value = " October 10 2020"
$ws.Range("C4").Value = "KKR won by 2 runs"
$ws.Range("D4").Value = "Kings XI Punjab"
$ws.Range("E4").Value = "Kolkata Knight Riders"
$ws.Range("F4").Value = "Prabhsimran Singh †"
$ws.Range("G4").Value = "'4"
$ws.Range("H4").Value = "'7"
$ws.Range("I4").Value = "'0"
$ws.Range("J4").Value = "'0"
$ws.Range("K4").Value = "'57.14"

$ws.Range("A5").Value = " Dubai (DSC)"
$ws.Range("B5").Value = " October 08 2020"
$ws.Range("C5").Value = "Sunrisers won by 69 runs"
$ws.Range("D5").Value = "Kings XI Punjab"
$ws.Range("E5").Value = "Sunrisers Hyderabad"
$ws.Range("F5").Value = "Prabhsimran Singh †"
$ws.Range("G5").Value = "'11"
$ws.Range("H5").Value = "'8"
$ws.Range("I5").Value = "'2"
$ws.Range("J5").Value = "'0"
$ws.Range("K5").Value = "'137.50"
